$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("D5").Copy()
$ws.Range("D4").PasteSpecial(-4122)
